# Reproduce the upload/re-save edit captured in the diff:
#   - the sole sheet is renamed from "C_29" to "C_38"
#   - the workbook window is scrolled back to the left edge (xWindow -> -120)
#
# (The dxfs re-shuffle and headerless xr:revisionPtr churn visible in the raw
# OOXML diff are artifacts Excel's own writer produces on every save/re-upload
# cycle - they carry no visible formatting change, so nothing below targets
# them directly.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the single worksheet tab.
$ws.Name = "C_38"

# Move the saved window position back towards the left of the screen.
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
